# Update hotel price figures on the "Hotel Data" sheet to reflect the
# latest rates reported by ExcelUtils.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hotel Data")

# Holiday Inn NAIROBI TWO RIVERS MALL by IHG
$ws.Range("C2").Value = "₹ 21,063"
$ws.Range("D2").Value = "₹ 134,802"

# JW Marriott Hotel Nairobi
$ws.Range("C3").Value = "₹ 206,210"
$ws.Range("D3").Value = "₹ 1,301,801"

# Yaya Hotel & Apartments
$ws.Range("C4").Value = "₹ 19,945"
$ws.Range("D4").Value = "₹ 96,167"
